$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.045256624090203
$ws.Range("D2").Value = 1.046813730444513
$ws.Range("E2").Value = 1.058612995940742
$ws.Range("F2").Value = 1.06554800064603
$ws.Range("I2").Value = 1.04154858850743
$ws.Range("J2").Value = 1.050317328946718
$ws.Range("K2").Value = 1.049578114640976
$ws.Range("L2").Value = 1.061344752243657
$ws.Range("M2").Value = 1.068260944868211
$ws.Range("N2").Value = 1.051808900026904
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.046247946365345
$ws.Range("D3").Value = 1.047567681846189
$ws.Range("E3").Value = 1.059626140512517
$ws.Range("F3").Value = 1.066702461307368
$ws.Range("I3").Value = 1.041796267816633
$ws.Range("J3").Value = 1.050956032849577
$ws.Range("K3").Value = 1.05014393438908
$ws.Range("L3").Value = 1.062171459722927
$ws.Range("M3").Value = 1.069229978555822
$ws.Range("N3").Value = 1.052448510962565
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04688969576673
$ws.Range("D4").Value = 1.048055714400943
$ws.Range("E4").Value = 1.060282415686127
$ws.Range("F4").Value = 1.067450485322019
$ws.Range("I4").Value = 1.0419554725832
$ws.Range("J4").Value = 1.05136899996446
$ws.Range("K4").Value = 1.050509568929103
$ws.Range("L4").Value = 1.062706484822704
$ws.Range("M4").Value = 1.069857433097031
$ws.Range("N4").Value = 1.052862064538119
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.04715955807125
$ws.Range("D5").Value = 1.048260924411808
$ws.Range("E5").Value = 1.060558481420158
$ws.Range("F5").Value = 1.067765196155145
$ws.Range("I5").Value = 1.042022148106349
$ws.Range("J5").Value = 1.051542534933536
$ws.Range("K5").Value = 1.050663164114864
$ws.Range("L5").Value = 1.062931430410937
$ws.Range("M5").Value = 1.070121316246855
$ws.Range("N5").Value = 1.053035845946756
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.047204873283448
$ws.Range("D6").Value = 1.048295382460743
$ws.Range("E6").Value = 1.060604843884801
$ws.Range("F6").Value = 1.067818051639942
$ws.Range("I6").Value = 1.042033328306691
$ws.Range("J6").Value = 1.051571667723153
$ws.Range("K6").Value = 1.050688946510774
$ws.Range("L6").Value = 1.062969200983289
$ws.Range("M6").Value = 1.070165629313197
$ws.Range("N6").Value = 1.053065020108277
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.046893301402109
$ws.Range("D7").Value = 1.048058456266698
$ws.Range("E7").Value = 1.060286103831992
$ws.Range("F7").Value = 1.067454689555039
$ws.Range("I7").Value = 1.041956364503877
$ws.Range("J7").Value = 1.051371319048068
$ws.Range("K7").Value = 1.050511621738212
$ws.Range("L7").Value = 1.06270949047466
$ws.Range("M7").Value = 1.069860958719291
$ws.Range("N7").Value = 1.052864386915092
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.045591584238212
$ws.Range("D8").Value = 1.04706849494808
$ws.Range("E8").Value = 1.058955246979747
$ws.Range("F8").Value = 1.065937946477279
$ws.Range("I8").Value = 1.041632512240036
$ws.Range("J8").Value = 1.05053324725423
$ws.Range("K8").Value = 1.0497694369043
$ws.Range("L8").Value = 1.06162412293622
$ws.Range("M8").Value = 1.068588346066343
$ws.Range("N8").Value = 1.052025124963178
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.043300088605116
$ws.Range("D9").Value = 1.04532544140459
$ws.Range("E9").Value = 1.056615517756782
$ws.Range("F9").Value = 1.06327301206961
$ws.Range("I9").Value = 1.04105373661432
$ws.Range("J9").Value = 1.049054052121931
$ws.Range("K9").Value = 1.048457895501874
$ws.Range("L9").Value = 1.059712279787817
$ws.Range("M9").Value = 1.066349115284687
$ws.Range("N9").Value = 1.050543829204125
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041773989578015
$ws.Range("D10").Value = 1.044164389510476
$ws.Range("E10").Value = 1.055059370950046
$ws.Range("F10").Value = 1.061501630812014
$ws.Range("I10").Value = 1.040662454188921
$ws.Range("J10").Value = 1.048066333549057
$ws.Range("K10").Value = 1.04758106907214
$ws.Range("L10").Value = 1.058438226760607
$ws.Range("M10").Value = 1.064858522830827
$ws.Range("N10").Value = 1.049554707957586
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.041113546650917
$ws.Range("D11").Value = 1.04366188529421
$ws.Range("E11").Value = 1.054386420047686
$ws.Range("F11").Value = 1.060735850697958
$ws.Range("I11").Value = 1.040491739818414
$ws.Range("J11").Value = 1.047638269367153
$ws.Range("K11").Value = 1.047200815444883
$ws.Range("L11").Value = 1.057886674159865
$ws.Range("M11").Value = 1.064213613120532
$ws.Range("N11").Value = 1.049126035875443
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040868284237976
$ws.Range("D12").Value = 1.043475269564692
$ws.Range("E12").Value = 1.054136587436701
$ws.Range("F12").Value = 1.060451592368502
$ws.Range("I12").Value = 1.040428135868493
$ws.Range("J12").Value = 1.047479211244053
$ws.Range("K12").Value = 1.047059485360902
$ws.Range("L12").Value = 1.057681821495865
$ws.Range("M12").Value = 1.063974144370618
$ws.Range("N12").Value = 1.048966751871564
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040920891343439
$ws.Range("D13").Value = 1.04351529759783
$ws.Range("E13").Value = 1.054190171429878
$ws.Range("F13").Value = 1.060512558301362
$ws.Range("I13").Value = 1.040441787862706
$ws.Range("J13").Value = 1.047513332280393
$ws.Range("K13").Value = 1.047089805070851
$ws.Range("L13").Value = 1.057725762242938
$ws.Range("M13").Value = 1.064025507630977
$ws.Range("N13").Value = 1.049000921363689
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.041093272050037
$ws.Range("D14").Value = 1.043646458826191
$ws.Range("E14").Value = 1.054365766124566
$ws.Range("F14").Value = 1.060712350009697
$ws.Range("I14").Value = 1.040486486232403
$ws.Range("J14").Value = 1.047625122713726
$ws.Range("K14").Value = 1.047189134827132
$ws.Range("L14").Value = 1.057869740594507
$ws.Range("M14").Value = 1.064193816936801
$ws.Range("N14").Value = 1.04911287055226
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.0411994888527
$ws.Range("D15").Value = 1.04372727647675
$ws.Range("E15").Value = 1.054473973226239
$ws.Range("F15").Value = 1.060835472996158
$ws.Range("I15").Value = 1.040514000802166
$ws.Range("J15").Value = 1.047693993076972
$ws.Range("K15").Value = 1.047250323660221
$ws.Range("L15").Value = 1.057958452875779
$ws.Range("M15").Value = 1.064297528389674
$ws.Range("N15").Value = 1.049181838719322
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041817828771172
$ws.Range("D16").Value = 1.04419774416368
$ws.Range("E16").Value = 1.055104050874677
$ws.Range("F16").Value = 1.061552479256115
$ws.Range("I16").Value = 1.040673756855712
$ws.Range("J16").Value = 1.04809473488056
$ws.Range("K16").Value = 1.047606293032174
$ws.Range("L16").Value = 1.058474834067507
$ws.Range("M16").Value = 1.064901334486494
$ws.Range("N16").Value = 1.049583149622237
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.042205795739118
$ws.Range("D17").Value = 1.044492920654705
$ws.Range("E17").Value = 1.05549951544116
$ws.Range("F17").Value = 1.06200257041154
$ws.Range("I17").Value = 1.040773623351147
$ws.Range("J17").Value = 1.048346009364312
$ws.Range("K17").Value = 1.047829427737571
$ws.Range("L17").Value = 1.058798778950752
$ws.Range("M17").Value = 1.065280227315286
$ws.Range("N17").Value = 1.049834780944577
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042432125802436
$ws.Range("D18").Value = 1.044665115105105
$ws.Range("E18").Value = 1.055730267232905
$ws.Range("F18").Value = 1.062265220559288
$ws.Range("I18").Value = 1.04083174964299
$ws.Range("J18").Value = 1.048492537260649
$ws.Range("K18").Value = 1.047959522249367
$ws.Range("L18").Value = 1.058987742146691
$ws.Range("M18").Value = 1.065501279886083
$ws.Range("N18").Value = 1.049981516927334
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042509304519496
$ws.Range("D19").Value = 1.04472383283868
$ws.Range("E19").Value = 1.0558089618645
$ws.Range("F19").Value = 1.062354797749288
$ws.Range("I19").Value = 1.040851548141615
$ws.Range("J19").Value = 1.048542493324681
$ws.Range("K19").Value = 1.048003871571051
$ws.Range("L19").Value = 1.059052175673585
$ws.Range("M19").Value = 1.065576661706515
$ws.Range("N19").Value = 1.050031543934707
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.042164166882183
$ws.Range("D20").Value = 1.044461248626096
$ws.Range("E20").Value = 1.055457077143754
$ws.Range("F20").Value = 1.061954267518311
$ws.Range("I20").Value = 1.040762921460419
$ws.Range("J20").Value = 1.048319053737782
$ws.Range("K20").Value = 1.047805493303535
$ws.Range("L20").Value = 1.058764021525155
$ws.Range("M20").Value = 1.065239570444398
$ws.Range("N20").Value = 1.049807787037966
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.041042508694773
$ws.Range("D21").Value = 1.043607834093606
$ws.Range("E21").Value = 1.054314054234775
$ws.Range("F21").Value = 1.060653511180836
$ws.Range("I21").Value = 1.040473329003202
$ws.Range("J21").Value = 1.047592204760857
$ws.Range("K21").Value = 1.047159887086164
$ws.Range("L21").Value = 1.057827342043906
$ws.Range("M21").Value = 1.064144251846062
$ws.Range("N21").Value = 1.049079905852123
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040337598864458
$ws.Range("D22").Value = 1.043071471304322
$ws.Range("E22").Value = 1.05359615036186
$ws.Range("F22").Value = 1.059836753214341
$ws.Range("I22").Value = 1.040290134035522
$ws.Range("J22").Value = 1.047134882004911
$ws.Range("K22").Value = 1.046753465972763
$ws.Range("L22").Value = 1.0572385218787
$ws.Range("M22").Value = 1.063456041289264
$ws.Range("N22").Value = 1.048621933645408
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040711254387806
$ws.Range("D23").Value = 1.043355786923816
$ws.Range("E23").Value = 1.053976652591557
$ws.Range("F23").Value = 1.060269629686126
$ws.Range("I23").Value = 1.040387354919737
$ws.Range("J23").Value = 1.047377348026516
$ws.Range("K23").Value = 1.046968964998513
$ws.Range("L23").Value = 1.05755065627896
$ws.Range("M23").Value = 1.063820831047177
$ws.Range("N23").Value = 1.048864743996575
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.042182977060983
$ws.Range("D24").Value = 1.044475559781503
$ws.Range("E24").Value = 1.055476252924032
$ws.Range("F24").Value = 1.061976093146985
$ws.Range("I24").Value = 1.040767757567762
$ws.Range("J24").Value = 1.04833123393713
$ws.Range("K24").Value = 1.047816308417635
$ws.Range("L24").Value = 1.058779726873703
$ws.Range("M24").Value = 1.065257941377287
$ws.Range("N24").Value = 1.049819984534594
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.043892220899341
$ws.Range("D25").Value = 1.045775892986519
$ws.Range("E25").Value = 1.05721974907675
$ws.Range("F25").Value = 1.063961039372963
$ws.Range("I25").Value = 1.041204322769693
$ws.Range("J25").Value = 1.049436741622925
$ws.Range("K25").Value = 1.048797397413198
$ws.Range("L25").Value = 1.060206449519709
$ws.Range("M25").Value = 1.066927619317483
$ws.Range("N25").Value = 1.050927062168106
